# Update DateBase/orders/International Ever Green_2024-12-24.xlsx
# - Append 10 new flower-order rows (82-91) to the "Orders" sheet (columns C/F).
# - Extend the sheet dimension / ignoredErrors range accordingly (A1:L81 -> A1:L91).
# - Append the corresponding digits to the "Number packed string" in Summary!G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)        # "Orders"
$summary = $wb.Worksheets.Item(2)   # "Summary"

# New rows to append at the bottom of the Orders table. Only columns C
# (FlowerName) and F (Number) are populated, same as every other line item
# row in the sheet.
$newRows = @(
    @{ Row = 82; C = "519_金鱼草粉色_snapdragon pink_undefined_1bunch"; F = "5" },
    @{ Row = 83; C = "7_翠绿洋桔梗_Dark Green Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "16" },
    @{ Row = 84; C = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "6" },
    @{ Row = 85; C = "13_酒红洋桔梗_Burgundy Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "10" },
    @{ Row = 86; C = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "7" },
    @{ Row = 87; C = "10_波浪粉洋桔梗_Wavy Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "5" },
    @{ Row = 88; C = "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "20" },
    @{ Row = 89; C = "579_腊梅红_wax red_undefined_1bunch"; F = "6" },
    @{ Row = 90; C = "577_腊梅白_wax white_undefined_1bunch"; F = "5" },
    @{ Row = 91; C = "578_腊梅粉_wax pink_undefined_1bunch"; F = $null }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # FlowerName values are never numeric-looking, so they come through as
    # text on their own.
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $item.C

    if ($item.F -ne $null) {
        # Number values ("5", "16", ...) look like real numbers, so force
        # text storage first (matches every other cell in this column,
        # which is stored as text even though it looks numeric) - otherwise
        # the COM layer would auto-coerce them into real numbers.
        $fCell = $ws.Cells.Item($r, 6)
        $fCell.NumberFormat = "@"
        $fCell.Value = $item.F
    }
}

# Summary!G2 packs every Orders!F value (blank treated as "0") back to back,
# prefixed with a leading "0". Extend it with the ten new rows' numbers.
$g2Cell = $summary.Range("G2")
$g2Cell.NumberFormat = "@"
$g2Cell.Value = "01013673102815383151019251841181010410197812530252525154060506101312251525655151210561010810553556581012515681251055203160151255166107520650"

Write-Output "Appended rows 82-91 to Orders and updated Summary!G2"
